# Auto-generated edit script: update Tonberry market-profit figures
# per-cell numeric updates (and a few cell deletions where the source
# sheet no longer has a computed profit figure), grouped by worksheet.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3591.8125
$ws.Range("J17").Value = 3201.4167
$ws.Range("L17").Value = 9604.250100000001
$ws.Range("N17").Value = -9940.250100000001
$ws.Range("H70").Value = 26620
$ws.Range("I70").Value = 950
$ws.Range("K70").Value = 2850
$ws.Range("M70").Value = -2580
$ws.Range("H73").Value = 26620
$ws.Range("I73").Value = 950
$ws.Range("K73").Value = 2850
$ws.Range("M73").Value = -1914
$ws.Range("H74").Value = 3737.375
$ws.Range("I74").Value = 2725
$ws.Range("K74").Value = 2725
$ws.Range("M74").Value = -1789
$ws.Range("H76").Value = 3500
$ws.Range("J76").Value = 3500
$ws.Range("L76").Value = 3500
$ws.Range("N76").Value = -4130
$ws.Range("H77").Value = 3737.375
$ws.Range("I77").Value = 2725
$ws.Range("K77").Value = 13625
$ws.Range("M77").Value = -8945
$ws.Range("H79").Value = 3500
$ws.Range("J79").Value = 3500
$ws.Range("L79").Value = 3500
$ws.Range("N79").Value = -5684
$ws.Range("H125").Value = 923
$ws.Range("I125").Value = 923
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 8307
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -5847
$ws.Range("N125").ClearContents()  # was -8304
$ws.Range("H135").Value = 394.35715
$ws.Range("I135").Value = 136.25
$ws.Range("K135").Value = 1226.25
$ws.Range("M135").Value = 1308.75
$ws.Range("H137").Value = 1263.5
$ws.Range("I137").Value = 921.4375
$ws.Range("K137").Value = 2764.3125
$ws.Range("M137").Value = -214.3125
$ws.Range("H138").Value = 2662.6604
$ws.Range("I138").Value = 2969.923
$ws.Range("J138").Value = 2366.7778
$ws.Range("K138").Value = 8909.769
$ws.Range("L138").Value = 7100.3334
$ws.Range("M138").Value = -3769.769
$ws.Range("N138").Value = -17380.3334

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1853184.6
$ws.Range("I2").Value = 5555555
$ws.Range("K2").Value = 5555555
$ws.Range("M2").Value = -5555442
$ws.Range("H32").Value = 5264.829
$ws.Range("I32").Value = 4317.154
$ws.Range("J32").Value = 23744.5
$ws.Range("K32").Value = 4317.154
$ws.Range("L32").Value = 23744.5
$ws.Range("M32").Value = -4030.154
$ws.Range("N32").Value = -24318.5
$ws.Range("H45").Value = 1665.2727
$ws.Range("I45").Value = 1041.6666
$ws.Range("K45").Value = 1041.6666
$ws.Range("M45").Value = -664.6666
$ws.Range("H74").Value = 1358.5
$ws.Range("I74").Value = 464.11765
$ws.Range("K74").Value = 464.11765
$ws.Range("M74").Value = 409.88235
$ws.Range("H77").Value = 1358.5
$ws.Range("I77").Value = 464.11765
$ws.Range("K77").Value = 2320.58825
$ws.Range("M77").Value = 2047.41175
$ws.Range("H97").Value = 1197.762
$ws.Range("I97").Value = 1119.7778
$ws.Range("K97").Value = 1119.7778
$ws.Range("M97").Value = -623.7778000000001
$ws.Range("H102").Value = 1511.5
$ws.Range("I102").Value = 1492
$ws.Range("J102").Value = 1521.25
$ws.Range("K102").Value = 1492
$ws.Range("L102").Value = 1521.25
$ws.Range("M102").Value = 130
$ws.Range("N102").Value = -4765.25
$ws.Range("H116").Value = 1853184.6
$ws.Range("I116").Value = 5555555
$ws.Range("K116").Value = 5555555
$ws.Range("M116").Value = -5553261
$ws.Range("H122").Value = 1231
$ws.Range("J122").Value = 1759.75
$ws.Range("L122").Value = 5279.25
$ws.Range("N122").Value = -10179.25
$ws.Range("H132").Value = 1745
$ws.Range("I132").Value = 1581.8
$ws.Range("K132").Value = 4745.4
$ws.Range("M132").Value = -2215.4

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1853184.6
$ws.Range("I3").Value = 5555555
$ws.Range("K3").Value = 5555555
$ws.Range("M3").Value = -5555441
$ws.Range("H99").Value = 1442.2222
$ws.Range("I99").Value = 1297
$ws.Range("J99").Value = 1623.75
$ws.Range("K99").Value = 1297
$ws.Range("L99").Value = 1623.75
$ws.Range("M99").Value = 201
$ws.Range("N99").Value = -4619.75
$ws.Range("H134").Value = 4585.5293
$ws.Range("I134").Value = 5050.815
$ws.Range("K134").Value = 15152.445
$ws.Range("M134").Value = -12617.445

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2229.4211
$ws.Range("I31").Value = 1915.8334
$ws.Range("K31").Value = 1915.8334
$ws.Range("M31").Value = -1620.8334
$ws.Range("H34").Value = 2229.4211
$ws.Range("I34").Value = 1915.8334
$ws.Range("K34").Value = 1915.8334
$ws.Range("M34").Value = -1713.8334
$ws.Range("H134").Value = 2212.182
$ws.Range("I134").Value = 2278.5557
$ws.Range("K134").Value = 6835.6671
$ws.Range("M134").Value = -4300.6671

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 379.72223
$ws.Range("I5").Value = 349.11765
$ws.Range("K5").Value = 1047.35295
$ws.Range("M5").Value = -935.35295
$ws.Range("H135").Value = 379.72223
$ws.Range("I135").Value = 349.11765
$ws.Range("K135").Value = 3142.05885
$ws.Range("M135").Value = -607.0588500000003

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 1801430
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10490
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()  # was -32994
$ws.Range("H102").Value = 1565.75
$ws.Range("I102").Value = 1238.88
$ws.Range("K102").Value = 1238.88
$ws.Range("M102").Value = 383.1199999999999
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()  # was -49820
$ws.Range("H126").Value = 2573264.8
$ws.Range("I126").Value = 3970959.5
$ws.Range("K126").Value = 11912878.5
$ws.Range("M126").Value = -11910408.5
$ws.Range("H132").Value = 2025465.2
$ws.Range("I132").Value = 2404933.8
$ws.Range("K132").Value = 7214801.399999999
$ws.Range("M132").Value = -7212271.399999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1634.5
$ws.Range("J46").Value = 2196.5
$ws.Range("L46").Value = 2196.5
$ws.Range("N46").Value = -2572.5
$ws.Range("H68").Value = 1474.6428
$ws.Range("I68").Value = 1203.4615
$ws.Range("K68").Value = 1203.4615
$ws.Range("M68").Value = -454.4614999999999
$ws.Range("H71").Value = 1474.6428
$ws.Range("I71").Value = 1203.4615
$ws.Range("K71").Value = 6017.307499999999
$ws.Range("M71").Value = -2273.307499999999
$ws.Range("H100").Value = 982.3333
$ws.Range("I100").Value = 982.3333
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 982.3333
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -441.3333
$ws.Range("N100").ClearContents()  # was -2114
$ws.Range("H122").Value = 9022.5
$ws.Range("I122").Value = 9108.076999999999
$ws.Range("K122").Value = 27324.231
$ws.Range("M122").Value = -24874.231
$ws.Range("H136").Value = 2384.8635
$ws.Range("I136").Value = 2304.923
$ws.Range("J136").Value = 2500.3333
$ws.Range("K136").Value = 6914.768999999999
$ws.Range("L136").Value = 7500.999899999999
$ws.Range("M136").Value = -4364.768999999999
$ws.Range("N136").Value = -12600.9999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()  # was -60663
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()  # was -376
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()  # was -1880
$ws.Range("H96").Value = 7918.9
$ws.Range("I96").Value = 372.25
$ws.Range("J96").Value = 12950
$ws.Range("K96").Value = 372.25
$ws.Range("L96").Value = 12950
$ws.Range("M96").Value = 1000.75
$ws.Range("N96").Value = -15696
$ws.Range("H113").Value = 495.95834
$ws.Range("I113").Value = 321.2143
$ws.Range("K113").Value = 963.6428999999999
$ws.Range("M113").Value = 1206.3571
$ws.Range("H136").Value = 1074.8975
$ws.Range("I136").Value = 694.4828
$ws.Range("J136").Value = 2178.1
$ws.Range("K136").Value = 2083.4484
$ws.Range("L136").Value = 6534.299999999999
$ws.Range("M136").Value = 466.5515999999998
$ws.Range("N136").Value = -11634.3

